$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.327452000000001
$ws.Range("H2").Value = 21.982356
$ws.Range("I2").Value = 0.2396705957528817
$ws.Range("J2").Value = 0.2396705957528817
$ws.Range("M2").Value = 10.32037433333333
$ws.Range("N2").Value = 30.961123
$ws.Range("O2").Value = 0.2929272563968436
$ws.Range("P2").Value = 0.2929272563968437
$ws.Range("Q2").Value = 75.62204754953201
$ws.Range("R2").Value = 680.5984279457881
$ws.Range("S2").Value = 0.07020605005288864
$ws.Range("T2").Value = 0.07020605005288866
$ws.Range("G3").Value = 7.327452000000001
$ws.Range("H3").Value = 21.982356
$ws.Range("I3").Value = 0.2396705957528817
$ws.Range("J3").Value = 0.2396705957528817
$ws.Range("M3").Value = 9.732885999999999
$ws.Range("O3").Value = 0.2762523432502674
$ws.Range("P3").Value = 0.2762523432502675
$ws.Range("Q3").Value = 71.31725498647199
$ws.Range("R3").Value = 641.8552948782481
$ws.Range("S3").Value = 0.06620956368492115
$ws.Range("T3").Value = 0.06620956368492117
$ws.Range("G4").Value = 7.327452000000001
$ws.Range("H4").Value = 21.982356
$ws.Range("I4").Value = 0.2396705957528817
$ws.Range("J4").Value = 0.2396705957528817
$ws.Range("M4").Value = 8.974543666666667
$ws.Range("N4").Value = 26.923631
$ws.Range("O4").Value = 0.2547280136147196
$ws.Range("P4").Value = 0.2547280136147196
$ws.Range("Q4").Value = 65.76053793940402
$ws.Range("R4").Value = 591.8448414546361
$ws.Range("S4").Value = 0.06105081477798801
$ws.Range("T4").Value = 0.06105081477798802
$ws.Range("G5").Value = 7.327452000000001
$ws.Range("H5").Value = 21.982356
$ws.Range("I5").Value = 0.2396705957528817
$ws.Range("J5").Value = 0.2396705957528817
$ws.Range("M5").Value = 6.204063666666666
$ws.Range("N5").Value = 18.612191
$ws.Range("O5").Value = 0.1760923867381692
$ws.Range("P5").Value = 0.1760923867381692
$ws.Range("Q5").Value = 45.45997872244401
$ws.Range("R5").Value = 409.139808501996
$ws.Range("S5").Value = 0.04220416723708385
$ws.Range("T5").Value = 0.04220416723708386
$ws.Range("I6").Value = 0.3687153231745697
$ws.Range("J6").Value = 0.3687153231745697
$ws.Range("M6").Value = 10.32037433333333
$ws.Range("N6").Value = 30.961123
$ws.Range("O6").Value = 0.2929272563968436
$ws.Range("P6").Value = 0.2929272563968437
$ws.Range("Q6").Value = 116.3388759215913
$ws.Range("R6").Value = 1047.049883294322
$ws.Range("S6").Value = 0.1080067680090022
$ws.Range("T6").Value = 0.1080067680090023
$ws.Range("I7").Value = 0.3687153231745697
$ws.Range("J7").Value = 0.3687153231745697
$ws.Range("M7").Value = 9.732885999999999
$ws.Range("O7").Value = 0.2762523432502674
$ws.Range("P7").Value = 0.2762523432502675
$ws.Range("R7").Value = 987.4464647568119
$ws.Range("S7").Value = 0.1018584720192545
$ws.Range("T7").Value = 0.1018584720192545
$ws.Range("I8").Value = 0.3687153231745697
$ws.Range("J8").Value = 0.3687153231745697
$ws.Range("M8").Value = 8.974543666666667
$ws.Range("N8").Value = 26.923631
$ws.Range("O8").Value = 0.2547280136147196
$ws.Range("P8").Value = 0.2547280136147196
$ws.Range("Q8").Value = 101.1676794238927
$ws.Range("R8").Value = 910.509114815034
$ws.Range("S8").Value = 0.09392212186156752
$ws.Range("T8").Value = 0.09392212186156755
$ws.Range("I9").Value = 0.3687153231745697
$ws.Range("J9").Value = 0.3687153231745697
$ws.Range("M9").Value = 6.204063666666666
$ws.Range("N9").Value = 18.612191
$ws.Range("O9").Value = 0.1760923867381692
$ws.Range("P9").Value = 0.1760923867381692
$ws.Range("Q9").Value = 69.93678424965266
$ws.Range("R9").Value = 629.4310582468739
$ws.Range("S9").Value = 0.06492796128474537
$ws.Range("T9").Value = 0.06492796128474539
$ws.Range("G10").Value = 6.386255666666667
$ws.Range("H10").Value = 19.158767
$ws.Range("I10").Value = 0.2088853943035337
$ws.Range("J10").Value = 0.2088853943035337
$ws.Range("M10").Value = 10.32037433333333
$ws.Range("N10").Value = 30.961123
$ws.Range("O10").Value = 0.2929272563968436
$ws.Range("P10").Value = 0.2929272563968437
$ws.Range("Q10").Value = 65.90854906837123
$ws.Range("R10").Value = 593.1769416153411
$ws.Range("S10").Value = 0.061188225454707
$ws.Range("T10").Value = 0.061188225454707
$ws.Range("G11").Value = 6.386255666666667
$ws.Range("H11").Value = 19.158767
$ws.Range("I11").Value = 0.2088853943035337
$ws.Range("J11").Value = 0.2088853943035337
$ws.Range("M11").Value = 9.732885999999999
$ws.Range("O11").Value = 0.2762523432502674
$ws.Range("P11").Value = 0.2762523432502675
$ws.Range("Q11").Value = 62.15669837052067
$ws.Range("R11").Value = 559.410285334686
$ws.Range("S11").Value = 0.05770507964710724
$ws.Range("T11").Value = 0.05770507964710725
$ws.Range("G12").Value = 6.386255666666667
$ws.Range("H12").Value = 19.158767
$ws.Range("I12").Value = 0.2088853943035337
$ws.Range("J12").Value = 0.2088853943035337
$ws.Range("M12").Value = 8.974543666666667
$ws.Range("N12").Value = 26.923631
$ws.Range("O12").Value = 0.2547280136147196
$ws.Range("P12").Value = 0.2547280136147196
$ws.Range("Q12").Value = 57.31373034699745
$ws.Range("R12").Value = 515.823573122977
$ws.Range("S12").Value = 0.0532089615640666
$ws.Range("T12").Value = 0.0532089615640666
$ws.Range("G13").Value = 6.386255666666667
$ws.Range("H13").Value = 19.158767
$ws.Range("I13").Value = 0.2088853943035337
$ws.Range("J13").Value = 0.2088853943035337
$ws.Range("M13").Value = 6.204063666666666
$ws.Range("N13").Value = 18.612191
$ws.Range("O13").Value = 0.1760923867381692
$ws.Range("P13").Value = 0.1760923867381692
$ws.Range("Q13").Value = 39.62073674761078
$ws.Range("R13").Value = 356.586630728497
$ws.Range("S13").Value = 0.03678312763765282
$ws.Range("T13").Value = 0.03678312763765283
$ws.Range("G14").Value = 5.586566333333333
$ws.Range("H14").Value = 16.759699
$ws.Range("I14").Value = 0.1827286867690149
$ws.Range("J14").Value = 0.1827286867690149
$ws.Range("M14").Value = 10.32037433333333
$ws.Range("N14").Value = 30.961123
$ws.Range("O14").Value = 0.2929272563968436
$ws.Range("P14").Value = 0.2929272563968437
$ws.Range("Q14").Value = 57.65545579799744
$ws.Range("R14").Value = 518.8991021819769
$ws.Range("S14").Value = 0.05352621288024575
$ws.Range("T14").Value = 0.05352621288024575
$ws.Range("G15").Value = 5.586566333333333
$ws.Range("H15").Value = 16.759699
$ws.Range("I15").Value = 0.1827286867690149
$ws.Range("J15").Value = 0.1827286867690149
$ws.Range("M15").Value = 9.732885999999999
$ws.Range("O15").Value = 0.2762523432502674
$ws.Range("P15").Value = 0.2762523432502675
$ws.Range("Q15").Value = 54.37341325377133
$ws.Range("R15").Value = 489.3607192839419
$ws.Range("S15").Value = 0.05047922789898449
$ws.Range("T15").Value = 0.0504792278989845
$ws.Range("G16").Value = 5.586566333333333
$ws.Range("H16").Value = 16.759699
$ws.Range("I16").Value = 0.1827286867690149
$ws.Range("J16").Value = 0.1827286867690149
$ws.Range("M16").Value = 8.974543666666667
$ws.Range("N16").Value = 26.923631
$ws.Range("O16").Value = 0.2547280136147196
$ws.Range("P16").Value = 0.2547280136147196
$ws.Range("Q16").Value = 50.13688350522989
$ws.Range("R16").Value = 451.2319515470689
$ws.Range("S16").Value = 0.04654611541109745
$ws.Range("T16").Value = 0.04654611541109745
$ws.Range("G17").Value = 5.586566333333333
$ws.Range("H17").Value = 16.759699
$ws.Range("I17").Value = 0.1827286867690149
$ws.Range("J17").Value = 0.1827286867690149
$ws.Range("M17").Value = 6.204063666666666
$ws.Range("N17").Value = 18.612191
$ws.Range("O17").Value = 0.1760923867381692
$ws.Range("P17").Value = 0.1760923867381692
$ws.Range("Q17").Value = 34.65941321005655
$ws.Range("R17").Value = 311.9347188905089
$ws.Range("S17").Value = 0.03217713057868715
$ws.Range("T17").Value = 0.03217713057868715
